# The sheet contains weekly price records for "Achicoria" (Femacal de La
# Calera). A new weekly record needs to be inserted as row 116, pushing the
# existing row 116 (and everything below it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 116, shifting rows 116:175 down to 117:176
$ws.Rows(116).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A116").Value = 3
$ws.Range("B116").Value = "Femacal de La Calera"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44572
$ws.Range("E116").Value = 5
$ws.Range("F116").Value = 100112010
$ws.Range("G116").Value = "Achicoria"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 115
$ws.Range("K116").Value = 5500
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 5739
$ws.Range("N116").Value = "$/caja 16 unidades"
$ws.Range("O116").Value = "Provincia de Quillota"
$ws.Range("P116").Value = 359
$ws.Range("Q116").Value = 16
$ws.Range("R116").Value = "Hortaliza"
